# Update Jenova_Profits market-data columns (H-N) across multiple item sheets
# Values come from a scheduled external data refresh; cells are plain numeric values (no formulas).
$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H5").Value = 132.58824
$ws_ALC.Range("I5").Value = 132.58824
$ws_ALC.Range("K5").Value = 132.58824
$ws_ALC.Range("M5").Value = -17.58824000000001

$ws_ALC.Range("H15").Value = 1673.5476
$ws_ALC.Range("I15").Value = 1673.5476
$ws_ALC.Range("K15").Value = 5020.642800000001
$ws_ALC.Range("M15").Value = -4851.642800000001

$ws_ALC.Range("H17").Value = 2324.7058
$ws_ALC.Range("J17").Value = 2324.7058
$ws_ALC.Range("L17").Value = 6974.117400000001
$ws_ALC.Range("N17").Value = -7310.117400000001

$ws_ALC.Range("H19").Value = 342.9091
$ws_ALC.Range("I19").Value = 143.66667
$ws_ALC.Range("K19").Value = 143.66667
$ws_ALC.Range("M19").Value = 31.33332999999999

$ws_ALC.Range("H64").Value = 8375
$ws_ALC.Range("I64").Value = 0
$ws_ALC.Range("K64").Value = 0
$ws_ALC.Range("M64").ClearContents()

$ws_ALC.Range("H67").Value = 8375
$ws_ALC.Range("I67").Value = 0
$ws_ALC.Range("K67").Value = 0
$ws_ALC.Range("M67").ClearContents()

$ws_ALC.Range("H98").Value = 2323.568
$ws_ALC.Range("I98").Value = 2100.2
$ws_ALC.Range("J98").Value = 3192.2222
$ws_ALC.Range("K98").Value = 2100.2
$ws_ALC.Range("L98").Value = 3192.2222
$ws_ALC.Range("M98").Value = -602.1999999999998
$ws_ALC.Range("N98").Value = -6188.2222

$ws_ALC.Range("H107").Value = 40700.56
$ws_ALC.Range("I107").Value = 40700.56
$ws_ALC.Range("J107").Value = 0
$ws_ALC.Range("K107").Value = 40700.56
$ws_ALC.Range("L107").Value = 0
$ws_ALC.Range("M107").Value = -38780.56
$ws_ALC.Range("N107").ClearContents()

$ws_ALC.Range("H113").Value = 4802.8335
$ws_ALC.Range("I113").Value = 4204.25
$ws_ALC.Range("K113").Value = 4204.25
$ws_ALC.Range("M113").Value = -950.25

$ws_ALC.Range("H115").Value = 444.3
$ws_ALC.Range("I115").Value = 444.3
$ws_ALC.Range("K115").Value = 1332.9
$ws_ALC.Range("M115").Value = 234.0999999999999

$ws_ALC.Range("H122").Value = 2323.568
$ws_ALC.Range("I122").Value = 2100.2
$ws_ALC.Range("J122").Value = 3192.2222
$ws_ALC.Range("K122").Value = 6300.599999999999
$ws_ALC.Range("L122").Value = 9576.6666
$ws_ALC.Range("M122").Value = -3850.599999999999
$ws_ALC.Range("N122").Value = -14476.6666

$ws_ALC.Range("H123").Value = 70772.664
$ws_ALC.Range("J123").Value = 70772.664
$ws_ALC.Range("L123").Value = 70772.664
$ws_ALC.Range("N123").Value = -80572.664

$ws_ALC.Range("H125").Value = 13893317
$ws_ALC.Range("J125").Value = 15877762
$ws_ALC.Range("L125").Value = 142899858
$ws_ALC.Range("N125").Value = -142904778

$ws_ALC.Range("H137").Value = 4290.409
$ws_ALC.Range("I137").Value = 4272.55
$ws_ALC.Range("J137").Value = 4469
$ws_ALC.Range("K137").Value = 12817.65
$ws_ALC.Range("L137").Value = 13407
$ws_ALC.Range("M137").Value = -10267.65
$ws_ALC.Range("N137").Value = -18507

$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H32").Value = 8295.5
$ws_ARM.Range("I32").Value = 8295.5
$ws_ARM.Range("K32").Value = 8295.5
$ws_ARM.Range("M32").Value = -8008.5

$ws_ARM.Range("I97").Value = 1549.375
$ws_ARM.Range("J97").Value = 1010
$ws_ARM.Range("K97").Value = 1549.375
$ws_ARM.Range("L97").Value = 1010
$ws_ARM.Range("M97").Value = -1053.375
$ws_ARM.Range("N97").Value = -2002

$ws_ARM.Range("H122").Value = 4759.025
$ws_ARM.Range("I122").Value = 4002.875
$ws_ARM.Range("K122").Value = 12008.625
$ws_ARM.Range("M122").Value = -9558.625

$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H86").Value = 1001774.7
$ws_BSM.Range("I86").Value = 1547251.9
$ws_BSM.Range("J86").Value = 1733.3334
$ws_BSM.Range("K86").Value = 1547251.9
$ws_BSM.Range("L86").Value = 1733.3334
$ws_BSM.Range("M86").Value = -1546128.9
$ws_BSM.Range("N86").Value = -3979.3334

$ws_BSM.Range("H89").Value = 1001774.7
$ws_BSM.Range("I89").Value = 1547251.9
$ws_BSM.Range("J89").Value = 1733.3334
$ws_BSM.Range("K89").Value = 7736259.5
$ws_BSM.Range("L89").Value = 8666.666999999999
$ws_BSM.Range("M89").Value = -7730643.5
$ws_BSM.Range("N89").Value = -19898.667

$ws_BSM.Range("H94").Value = 132.5
$ws_BSM.Range("I94").Value = 255
$ws_BSM.Range("J94").Value = 10
$ws_BSM.Range("K94").Value = 255
$ws_BSM.Range("L94").Value = 10
$ws_BSM.Range("M94").Value = 196
$ws_BSM.Range("N94").Value = -912

$ws_BSM.Range("H105").Value = 59703.65
$ws_BSM.Range("J105").Value = 1218
$ws_BSM.Range("L105").Value = 1218
$ws_BSM.Range("N105").Value = -4712

$ws_BSM.Range("H134").Value = 25294.848
$ws_BSM.Range("I134").Value = 3863.8918
$ws_BSM.Range("K134").Value = 11591.6754
$ws_BSM.Range("M134").Value = -9056.6754

$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H62").Value = 4398
$ws_CRP.Range("I62").Value = 3497.5
$ws_CRP.Range("J62").Value = 8000
$ws_CRP.Range("K62").Value = 3497.5
$ws_CRP.Range("L62").Value = 8000
$ws_CRP.Range("M62").Value = -2873.5
$ws_CRP.Range("N62").Value = -9248

$ws_CRP.Range("H65").Value = 4398
$ws_CRP.Range("I65").Value = 3497.5
$ws_CRP.Range("J65").Value = 8000
$ws_CRP.Range("K65").Value = 17487.5
$ws_CRP.Range("L65").Value = 40000
$ws_CRP.Range("M65").Value = -14367.5
$ws_CRP.Range("N65").Value = -46240

$ws_CRP.Range("H122").Value = 3788.818
$ws_CRP.Range("J122").Value = 4538.6
$ws_CRP.Range("L122").Value = 13615.8
$ws_CRP.Range("N122").Value = -18515.8

$ws_CRP.Range("H132").Value = 1386.6
$ws_CRP.Range("I132").Value = 1283.6
$ws_CRP.Range("J132").Value = 1695.6
$ws_CRP.Range("K132").Value = 3850.8
$ws_CRP.Range("L132").Value = 5086.799999999999
$ws_CRP.Range("M132").Value = -1320.8
$ws_CRP.Range("N132").Value = -10146.8

$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H131").Value = 3792.2683
$ws_CUL.Range("I131").Value = 1206.75
$ws_CUL.Range("J131").Value = 4419.0605
$ws_CUL.Range("K131").Value = 3620.25
$ws_CUL.Range("L131").Value = 13257.1815
$ws_CUL.Range("M131").Value = 1419.75
$ws_CUL.Range("N131").Value = -23337.1815

$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H99").Value = 4471
$ws_GSM.Range("I99").Value = 4471
$ws_GSM.Range("K99").Value = 4471
$ws_GSM.Range("M99").Value = -2225

$ws_GSM.Range("H132").Value = 52725.76
$ws_GSM.Range("I132").Value = 5809.2354
$ws_GSM.Range("K132").Value = 17427.7062
$ws_GSM.Range("M132").Value = -14897.7062

$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H68").Value = 144807.72
$ws_LTW.Range("I68").Value = 1199.6666
$ws_LTW.Range("K68").Value = 1199.6666
$ws_LTW.Range("M68").Value = -450.6666

$ws_LTW.Range("H71").Value = 144807.72
$ws_LTW.Range("I71").Value = 1199.6666
$ws_LTW.Range("K71").Value = 5998.333000000001
$ws_LTW.Range("M71").Value = -2254.333000000001

$ws_LTW.Range("H122").Value = 1054930
$ws_LTW.Range("I122").Value = 557981.4
$ws_LTW.Range("K122").Value = 1673944.2
$ws_LTW.Range("M122").Value = -1671494.2

$ws_LTW.Range("H123").Value = 80306.336
$ws_LTW.Range("J123").Value = 80306.336
$ws_LTW.Range("L123").Value = 80306.336
$ws_LTW.Range("N123").Value = -90106.336

$ws_LTW.Range("H124").Value = 67161.60000000001
$ws_LTW.Range("J124").Value = 67161.60000000001
$ws_LTW.Range("L124").Value = 67161.60000000001
$ws_LTW.Range("N124").Value = -76981.60000000001

$ws_LTW.Range("H125").Value = 81250
$ws_LTW.Range("J125").Value = 81250
$ws_LTW.Range("L125").Value = 81250
$ws_LTW.Range("N125").Value = -91090

$ws_LTW.Range("H132").Value = 6506.696
$ws_LTW.Range("I132").Value = 5986.625
$ws_LTW.Range("K132").Value = 17959.875
$ws_LTW.Range("M132").Value = -15429.875

$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H42").Value = 0
$ws_WVR.Range("I42").Value = 0
$ws_WVR.Range("K42").Value = 0
$ws_WVR.Range("M42").ClearContents()

$ws_WVR.Range("H81").Value = 1661.4667
$ws_WVR.Range("I81").Value = 1385.8334
$ws_WVR.Range("J81").Value = 2764
$ws_WVR.Range("K81").Value = 2771.6668
$ws_WVR.Range("L81").Value = 5528
$ws_WVR.Range("M81").Value = -1710.6668
$ws_WVR.Range("N81").Value = -7650

$ws_WVR.Range("H84").Value = 1661.4667
$ws_WVR.Range("I84").Value = 1385.8334
$ws_WVR.Range("J84").Value = 2764
$ws_WVR.Range("K84").Value = 13858.334
$ws_WVR.Range("L84").Value = 27640
$ws_WVR.Range("M84").Value = -8554.333999999999
$ws_WVR.Range("N84").Value = -38248

$ws_WVR.Range("H88").Value = 24414.285
$ws_WVR.Range("J88").Value = 28786.2
$ws_WVR.Range("L88").Value = 28786.2
$ws_WVR.Range("N88").Value = -29598.2

$ws_WVR.Range("H91").Value = 24414.285
$ws_WVR.Range("J91").Value = 28786.2
$ws_WVR.Range("L91").Value = 28786.2
$ws_WVR.Range("N91").Value = -31594.2

$ws_WVR.Range("H113").Value = 388.1875
$ws_WVR.Range("I113").Value = 393.23077
$ws_WVR.Range("J113").Value = 366.33334
$ws_WVR.Range("K113").Value = 1179.69231
$ws_WVR.Range("L113").Value = 1099.00002
$ws_WVR.Range("M113").Value = 990.3076900000001
$ws_WVR.Range("N113").Value = -5439.000019999999

Write-Output "Jenova_Profits sheets updated"